$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (columns D..M = col index 4..13) ---
$ws.Cells.Item(8, 4).Value = "12 ماهه منتهی به 1399/09"
$ws.Cells.Item(8, 5).Value = "3 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 6).Value = "6 ماهه منتهی به 1400/03"
$ws.Cells.Item(8, 7).Value = "9 ماهه منتهی به 1400/06"
$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1400/09"
$ws.Cells.Item(8, 9).Value = "3 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 10).Value = "6 ماهه منتهی به 1401/03"
$ws.Cells.Item(8, 11).Value = "9 ماهه منتهی به 1401/06"
$ws.Cells.Item(8, 12).Value = "12 ماهه منتهی به 1401/09"
$ws.Cells.Item(8, 13).Value = "3 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (columns D..M) ---
$ws.Cells.Item(9, 4).Value = "1400-11-27 (11)"
$ws.Cells.Item(9, 5).Value = "1401-01-30 (2)"
$ws.Cells.Item(9, 6).Value = "1401-06-15 (6)"
$ws.Cells.Item(9, 7).Value = "1401-08-06 (3)"
$ws.Cells.Item(9, 8).Value = "1401-12-16 (8)"
$ws.Cells.Item(9, 9).Value = "1402-02-03 (2)"
$ws.Cells.Item(9, 10).Value = "1401-06-15 (2)"
$ws.Cells.Item(9, 11).Value = "1401-08-06"
$ws.Cells.Item(9, 12).Value = "1402-02-03 (3)"
$ws.Cells.Item(9, 13).Value = "1402-02-03"

# --- Data rows 11-27: shift quarterly columns left and append new quarter (columns D..M) ---
# Row 11
$ws.Cells.Item(11, 4).Value = 8195863
$ws.Cells.Item(11, 5).Value = 2469916
$ws.Cells.Item(11, 6).Value = 5018726
$ws.Cells.Item(11, 7).Value = 8503541
$ws.Cells.Item(11, 8).Value = 11206859
$ws.Cells.Item(11, 9).Value = 2319378
$ws.Cells.Item(11, 10).Value = 3485694
$ws.Cells.Item(11, 11).Value = 6704606
$ws.Cells.Item(11, 12).Value = 11234217
$ws.Cells.Item(11, 13).Value = 4762946

# Row 12
$ws.Cells.Item(12, 4).Value = -4851597
$ws.Cells.Item(12, 5).Value = -1695770
$ws.Cells.Item(12, 6).Value = -3619987
$ws.Cells.Item(12, 7).Value = -5974089
$ws.Cells.Item(12, 8).Value = -8023476
$ws.Cells.Item(12, 9).Value = -2040833
$ws.Cells.Item(12, 10).Value = -2912462
$ws.Cells.Item(12, 11).Value = -5165847
$ws.Cells.Item(12, 12).Value = -8366153
$ws.Cells.Item(12, 13).Value = -3732414

# Row 13
$ws.Cells.Item(13, 4).Value = 3344266
$ws.Cells.Item(13, 5).Value = 774146
$ws.Cells.Item(13, 6).Value = 1398739
$ws.Cells.Item(13, 7).Value = 2529452
$ws.Cells.Item(13, 8).Value = 3183383
$ws.Cells.Item(13, 9).Value = 278545
$ws.Cells.Item(13, 10).Value = 573232
$ws.Cells.Item(13, 11).Value = 1538759
$ws.Cells.Item(13, 12).Value = 2868064
$ws.Cells.Item(13, 13).Value = 1030532

# Row 14
$ws.Cells.Item(14, 4).Value = -814284
$ws.Cells.Item(14, 5).Value = -340351
$ws.Cells.Item(14, 6).Value = -647878
$ws.Cells.Item(14, 7).Value = -970850
$ws.Cells.Item(14, 8).Value = -1414297
$ws.Cells.Item(14, 9).Value = -244782
$ws.Cells.Item(14, 10).Value = -401953
$ws.Cells.Item(14, 11).Value = -612489
$ws.Cells.Item(14, 12).Value = -1038150
$ws.Cells.Item(14, 13).Value = -352062

# Row 15
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0

# Row 16
$ws.Cells.Item(16, 4).Value = 26270
$ws.Cells.Item(16, 5).Value = 17369
$ws.Cells.Item(16, 6).Value = 18035
$ws.Cells.Item(16, 7).Value = 4161
$ws.Cells.Item(16, 8).Value = 3164
$ws.Cells.Item(16, 9).Value = 890
$ws.Cells.Item(16, 10).Value = -357
$ws.Cells.Item(16, 11).Value = 1967
$ws.Cells.Item(16, 12).Value = 8721
$ws.Cells.Item(16, 13).Value = -24655

# Row 17
$ws.Cells.Item(17, 4).Value = 2556252
$ws.Cells.Item(17, 5).Value = 451164
$ws.Cells.Item(17, 6).Value = 768896
$ws.Cells.Item(17, 7).Value = 1562763
$ws.Cells.Item(17, 8).Value = 1772250
$ws.Cells.Item(17, 9).Value = 34653
$ws.Cells.Item(17, 10).Value = 170922
$ws.Cells.Item(17, 11).Value = 928237
$ws.Cells.Item(17, 12).Value = 1838635
$ws.Cells.Item(17, 13).Value = 653815

# Row 18
$ws.Cells.Item(18, 4).Value = -117381
$ws.Cells.Item(18, 5).Value = -90348
$ws.Cells.Item(18, 6).Value = -245230
$ws.Cells.Item(18, 7).Value = -429880
$ws.Cells.Item(18, 8).Value = -633422
$ws.Cells.Item(18, 9).Value = -211262
$ws.Cells.Item(18, 10).Value = -401495
$ws.Cells.Item(18, 11).Value = -542861
$ws.Cells.Item(18, 12).Value = -911067
$ws.Cells.Item(18, 13).Value = -113077

# Row 19
$ws.Cells.Item(19, 4).Value = 57272
$ws.Cells.Item(19, 5).Value = 7967
$ws.Cells.Item(19, 6).Value = 91806
$ws.Cells.Item(19, 7).Value = 96790
$ws.Cells.Item(19, 8).Value = 43927
$ws.Cells.Item(19, 9).Value = 3355
$ws.Cells.Item(19, 10).Value = 4491
$ws.Cells.Item(19, 11).Value = 9907
$ws.Cells.Item(19, 12).Value = 17654
$ws.Cells.Item(19, 13).Value = 3604

# Row 20
$ws.Cells.Item(20, 4).Value = 2496143
$ws.Cells.Item(20, 5).Value = 368783
$ws.Cells.Item(20, 6).Value = 615472
$ws.Cells.Item(20, 7).Value = 1229673
$ws.Cells.Item(20, 8).Value = 1182755
$ws.Cells.Item(20, 9).Value = -173254
$ws.Cells.Item(20, 10).Value = -226082
$ws.Cells.Item(20, 11).Value = 395283
$ws.Cells.Item(20, 12).Value = 945222
$ws.Cells.Item(20, 13).Value = 544342

# Row 21
$ws.Cells.Item(21, 4).Value = -440467
$ws.Cells.Item(21, 5).Value = -71456
$ws.Cells.Item(21, 6).Value = -152966
$ws.Cells.Item(21, 7).Value = -170230
$ws.Cells.Item(21, 8).Value = -199416
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = -100536
$ws.Cells.Item(21, 11).Value = -172258
$ws.Cells.Item(21, 12).Value = -241333
$ws.Cells.Item(21, 13).Value = -70606

# Row 22
$ws.Cells.Item(22, 4).Value = 2055676
$ws.Cells.Item(22, 5).Value = 297327
$ws.Cells.Item(22, 6).Value = 462506
$ws.Cells.Item(22, 7).Value = 1059443
$ws.Cells.Item(22, 8).Value = 983339
$ws.Cells.Item(22, 9).Value = -173254
$ws.Cells.Item(22, 10).Value = -326618
$ws.Cells.Item(22, 11).Value = 223025
$ws.Cells.Item(22, 12).Value = 703889
$ws.Cells.Item(22, 13).Value = 473736

# Row 23
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = 0

# Row 24
$ws.Cells.Item(24, 4).Value = 2055676
$ws.Cells.Item(24, 5).Value = 297327
$ws.Cells.Item(24, 6).Value = 462506
$ws.Cells.Item(24, 7).Value = 1059443
$ws.Cells.Item(24, 8).Value = 983339
$ws.Cells.Item(24, 9).Value = -173254
$ws.Cells.Item(24, 10).Value = -326618
$ws.Cells.Item(24, 11).Value = 223025
$ws.Cells.Item(24, 12).Value = 703889
$ws.Cells.Item(24, 13).Value = 473736

# Row 25
$ws.Cells.Item(25, 4).Value = 791
$ws.Cells.Item(25, 5).Value = 114
$ws.Cells.Item(25, 6).Value = 183
$ws.Cells.Item(25, 7).Value = 407
$ws.Cells.Item(25, 8).Value = 378
$ws.Cells.Item(25, 9).Value = -67
$ws.Cells.Item(25, 10).Value = -126
$ws.Cells.Item(25, 11).Value = 86
$ws.Cells.Item(25, 12).Value = 271
$ws.Cells.Item(25, 13).Value = 182

# Row 26
$ws.Cells.Item(26, 4).Value = 2600000
$ws.Cells.Item(26, 5).Value = 2600000
$ws.Cells.Item(26, 6).Value = 2523560
$ws.Cells.Item(26, 7).Value = 2600000
$ws.Cells.Item(26, 8).Value = 2600000
$ws.Cells.Item(26, 9).Value = 2600000
$ws.Cells.Item(26, 10).Value = 2600000
$ws.Cells.Item(26, 11).Value = 2600000
$ws.Cells.Item(26, 12).Value = 2600000
$ws.Cells.Item(26, 13).Value = 2600000

# Row 27
$ws.Cells.Item(27, 4).Value = 791
$ws.Cells.Item(27, 5).Value = 114
$ws.Cells.Item(27, 6).Value = 178
$ws.Cells.Item(27, 7).Value = 407
$ws.Cells.Item(27, 8).Value = 378
$ws.Cells.Item(27, 9).Value = -67
$ws.Cells.Item(27, 10).Value = -126
$ws.Cells.Item(27, 11).Value = 86
$ws.Cells.Item(27, 12).Value = 271
$ws.Cells.Item(27, 13).Value = 182

# --- Column widths: shift pattern matches the new header/date text widths ---
$ws.Columns.Item(4).ColumnWidth = 28.166666666666668
$ws.Columns.Item(5).ColumnWidth = 27.166666666666668
$ws.Columns.Item(6).ColumnWidth = 27.166666666666668
$ws.Columns.Item(7).ColumnWidth = 27.166666666666668
$ws.Columns.Item(8).ColumnWidth = 28.166666666666668
$ws.Columns.Item(9).ColumnWidth = 27.166666666666668
$ws.Columns.Item(10).ColumnWidth = 27.166666666666668
$ws.Columns.Item(11).ColumnWidth = 27.166666666666668
$ws.Columns.Item(12).ColumnWidth = 28.166666666666668
$ws.Columns.Item(13).ColumnWidth = 27.166666666666668
